# Add MinD / MaxD columns (Min/Max Depth) to the pings table on sheet "Data".
# MinD = 0.8 * Mean Depth (m), MaxD = 1.2 * Mean Depth (m), rounded to 1 decimal.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Header row (row 8): copy the style used by the other headers (e.g. F8) and
# set the new header text.
$ws.Range("G8").Value = "MinD"
$ws.Range("H8").Value = "MaxD"
$ws.Range("G8:H8").Font.Bold = $true
$ws.Range("G8:H8").Borders.LineStyle = 1

# Data rows 9-23: Min/Max depth (m), derived from the underlying (unrounded)
# ping-depth statistics -- roughly Mean Depth +/- 20%, rounded to 1 decimal.
$minMax = @{
    9  = @(5.3, 7.9)
    10 = @(24.5, 36.7)
    11 = @(138.2, 207.4)
    12 = @(93.8, 140.6)
    13 = @(117.8, 176.6)
    14 = @(101, 151.4)
    15 = @(96, 143.9)
    16 = @(124.7, 187.1)
    17 = @(120.4, 180.6)
    18 = @(26, 39)
    19 = @(80.6, 121)
    20 = @(116.6, 174.8)
    21 = @(93.8, 140.8)
    22 = @(92.4, 138.6)
    23 = @(1.7, 2.5)
}

for ($r = 9; $r -le 23; $r++) {
    $ws.Range("G$r").Value = $minMax[$r][0]
    $ws.Range("H$r").Value = $minMax[$r][1]
}
$ws.Range("G9:H23").Borders.LineStyle = 1

$ws.Range("H4").Select()
